# Update the Lgi1-Adam11 NATMI worksheet with new TPM-derived values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-obsolete rows (previously rows 5-7: MuSCs, Neutrophils, Resolving-Mac
# target-cluster entries). Deleting shifts rows 8+ up, but there are none here, so the
# sheet simply shrinks from 7 data+header rows down to 4.
$ws.Range("A5:T7").Delete()

# Row 2 (Sending cluster "Neutrophils" -> "MuSCs"; Target cluster stays "ECs")
$ws.Range("A2").Value = "MuSCs"
$ws.Range("F2").Value = 0.5
$ws.Range("G2").Value = 0.035285
$ws.Range("H2").Value = 0.07056999999999999
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.4578845
$ws.Range("N2").Value = 0.9157690000000001
$ws.Range("O2").Value = 0.1062139753234554
$ws.Range("P2").Value = 0.1055270602873987
$ws.Range("Q2").Value = 0.0161564545825
$ws.Range("R2").Value = 0.06462581833
$ws.Range("S2").Value = 0.1062139753234554
$ws.Range("T2").Value = 0.1055270602873987

# Row 3 (Sending cluster "Neutrophils" -> "MuSCs"; Target cluster stays "FAPs")
$ws.Range("A3").Value = "MuSCs"
$ws.Range("F3").Value = 0.5
$ws.Range("G3").Value = 0.035285
$ws.Range("H3").Value = 0.07056999999999999
$ws.Range("O3").Value = 0.01301874674014239
$ws.Range("P3").Value = 0.01940182637825622
$ws.Range("Q3").Value = 0.001980311816666666
$ws.Range("R3").Value = 0.0118818709
$ws.Range("S3").Value = 0.01301874674014239
$ws.Range("T3").Value = 0.01940182637825622

# Row 4 (Sending cluster "Neutrophils" -> "MuSCs"; Target cluster "Inflammatory-Mac" -> "MuSCs")
$ws.Range("A4").Value = "MuSCs"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("F4").Value = 0.5
$ws.Range("G4").Value = 0.035285
$ws.Range("H4").Value = 0.07056999999999999
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.796955
$ws.Range("N4").Value = 7.593909999999999
$ws.Range("O4").Value = 0.8807672779364022
$ws.Range("P4").Value = 0.8750711133343451
$ws.Range("Q4").Value = 0.133975557175
$ws.Range("R4").Value = 0.5359022286999999
$ws.Range("S4").Value = 0.8807672779364022
$ws.Range("T4").Value = 0.8750711133343451
